$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$tb = $s.Shapes.Item(3)
$tf = $tb.TextFrame
$tr = $tf.TextRange

# --- Paragraph 2 ("Создать базу данных ... по ссылке (...):  https://...") ---
$para2 = $tr.Paragraphs(2)
$base2 = $para2.Start

# 1) Merge "Создать " + "базу данных" -> "Создать базу данных"
$tr.Characters($base2 + 0, 19).Text = "Создать базу данных"

# 2) Merge " " + "по ссылке (" -> " по ссылке ("
$tr.Characters($base2 + 54, 12).Text = " по ссылке ("

# 3) Merge "): " + " " -> "):  " (two trailing spaces)
$tr.Characters($base2 + 107, 4).Text = "):  "

# --- Paragraph 4 ("Выполнить работу до 23.09.2022г. ...") ---
# Re-fetch paragraph 4 start (unaffected by the same-length edits above, but re-query to be safe)
$para4 = $tr.Paragraphs(4)
$base4 = $para4.Start

$tr.Characters($base4 + 0, 17).Text = "Выполнить работу "
$tr.Characters($base4 + 17, 3).Text = "до "
$tr.Characters($base4 + 20, 11).Text = "26.09.2022г"
$tr.Characters($base4 + 31, 25).Text = ". Включительно. Ссылку в "
